$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 51.56869633333334
$ws.Cells.Item(2, 8).Value = 154.706089
$ws.Cells.Item(2, 9).Value = 0.1855839901450455
$ws.Cells.Item(2, 10).Value = 0.1855839901450455
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 1.307106666666667
$ws.Cells.Item(2, 14).Value = 3.92132
$ws.Cells.Item(2, 15).Value = 0.01256263154946851
$ws.Cells.Item(2, 16).Value = 0.01256263154946851
$ws.Cells.Item(2, 17).Value = 67.40578676860891
$ws.Cells.Item(2, 18).Value = 606.6520809174801
$ws.Cells.Item(2, 19).Value = 0.002331423289672401
$ws.Cells.Item(2, 20).Value = 0.002331423289672401

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 51.56869633333334
$ws.Cells.Item(3, 8).Value = 154.706089
$ws.Cells.Item(3, 9).Value = 0.1855839901450455
$ws.Cells.Item(3, 10).Value = 0.1855839901450455
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 80.22623699999998
$ws.Cells.Item(3, 14).Value = 240.678711
$ws.Cells.Item(3, 15).Value = 0.77105616682495
$ws.Cells.Item(3, 16).Value = 0.77105616682495
$ws.Cells.Item(3, 17).Value = 4137.162453819031
$ws.Cells.Item(3, 18).Value = 37234.46208437128
$ws.Cells.Item(3, 19).Value = 0.1430956800653181
$ws.Cells.Item(3, 20).Value = 0.1430956800653181

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 51.56869633333334
$ws.Cells.Item(4, 8).Value = 154.706089
$ws.Cells.Item(4, 9).Value = 0.1855839901450455
$ws.Cells.Item(4, 10).Value = 0.1855839901450455
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 22.51385866666667
$ws.Cells.Item(4, 14).Value = 67.54157600000001
$ws.Cells.Item(4, 15).Value = 0.2163812016255815
$ws.Cells.Item(4, 16).Value = 0.2163812016255815
$ws.Cells.Item(4, 17).Value = 1161.010340872919
$ws.Cells.Item(4, 18).Value = 10449.09306785627
$ws.Cells.Item(4, 19).Value = 0.040156886790055
$ws.Cells.Item(4, 20).Value = 0.040156886790055

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 165.03522
$ws.Cells.Item(5, 8).Value = 495.1056600000001
$ws.Cells.Item(5, 9).Value = 0.5939241598059933
$ws.Cells.Item(5, 10).Value = 0.5939241598059933
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.307106666666667
$ws.Cells.Item(5, 14).Value = 3.92132
$ws.Cells.Item(5, 15).Value = 0.01256263154946851
$ws.Cells.Item(5, 16).Value = 0.01256263154946851
$ws.Cells.Item(5, 17).Value = 215.7186362968
$ws.Cells.Item(5, 18).Value = 1941.4677266712
$ws.Cells.Item(5, 19).Value = 0.007461250387970349
$ws.Cells.Item(5, 20).Value = 0.007461250387970349

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 165.03522
$ws.Cells.Item(6, 8).Value = 495.1056600000001
$ws.Cells.Item(6, 9).Value = 0.5939241598059933
$ws.Cells.Item(6, 10).Value = 0.5939241598059933
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 80.22623699999998
$ws.Cells.Item(6, 14).Value = 240.678711
$ws.Cells.Item(6, 15).Value = 0.77105616682495
$ws.Cells.Item(6, 16).Value = 0.77105616682495
$ws.Cells.Item(6, 17).Value = 13240.15467306714
$ws.Cells.Item(6, 18).Value = 119161.3920576043
$ws.Cells.Item(6, 19).Value = 0.4579488860447383
$ws.Cells.Item(6, 20).Value = 0.4579488860447383

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 165.03522
$ws.Cells.Item(7, 8).Value = 495.1056600000001
$ws.Cells.Item(7, 9).Value = 0.5939241598059933
$ws.Cells.Item(7, 10).Value = 0.5939241598059933
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 22.51385866666667
$ws.Cells.Item(7, 14).Value = 67.54157600000001
$ws.Cells.Item(7, 15).Value = 0.2163812016255815
$ws.Cells.Item(7, 16).Value = 0.2163812016255815
$ws.Cells.Item(7, 17).Value = 3715.57961810224
$ws.Cells.Item(7, 18).Value = 33440.21656292017
$ws.Cells.Item(7, 19).Value = 0.1285140233732847
$ws.Cells.Item(7, 20).Value = 0.1285140233732847

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 61.26863233333334
$ws.Cells.Item(8, 8).Value = 183.805897
$ws.Cells.Item(8, 9).Value = 0.2204918500489612
$ws.Cells.Item(8, 10).Value = 0.2204918500489612
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.307106666666667
$ws.Cells.Item(8, 14).Value = 3.92132
$ws.Cells.Item(8, 15).Value = 0.01256263154946851
$ws.Cells.Item(8, 16).Value = 0.01256263154946851
$ws.Cells.Item(8, 17).Value = 80.08463778044889
$ws.Cells.Item(8, 18).Value = 720.7617400240401
$ws.Cells.Item(8, 19).Value = 0.002769957871825759
$ws.Cells.Item(8, 20).Value = 0.002769957871825759

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 61.26863233333334
$ws.Cells.Item(9, 8).Value = 183.805897
$ws.Cells.Item(9, 9).Value = 0.2204918500489612
$ws.Cells.Item(9, 10).Value = 0.2204918500489612
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 80.22623699999998
$ws.Cells.Item(9, 14).Value = 240.678711
$ws.Cells.Item(9, 15).Value = 0.77105616682495
$ws.Cells.Item(9, 16).Value = 0.77105616682495
$ws.Cells.Item(9, 17).Value = 4915.351818239862
$ws.Cells.Item(9, 18).Value = 44238.16636415877
$ws.Cells.Item(9, 19).Value = 0.1700116007148937
$ws.Cells.Item(9, 20).Value = 0.1700116007148937

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 61.26863233333334
$ws.Cells.Item(10, 8).Value = 183.805897
$ws.Cells.Item(10, 9).Value = 0.2204918500489612
$ws.Cells.Item(10, 10).Value = 0.2204918500489612
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 22.51385866666667
$ws.Cells.Item(10, 14).Value = 67.54157600000001
$ws.Cells.Item(10, 15).Value = 0.2163812016255815
$ws.Cells.Item(10, 16).Value = 0.2163812016255815
$ws.Cells.Item(10, 17).Value = 1379.39332905263
$ws.Cells.Item(10, 18).Value = 12414.53996147367
$ws.Cells.Item(10, 19).Value = 0.04771029146224174
$ws.Cells.Item(10, 20).Value = 0.04771029146224174
